$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "25 x 51" + [char]11 + "  5    1" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "5|    |"
$t.Cell(1,2).Range.Text = "60 x 81" + [char]11 + "  8    1" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "0|    |"
$t.Cell(1,3).Range.Text = "83 x 91" + [char]11 + "  9    1" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "3|    |"
$t.Cell(2,1).Range.Text = "58 x 34" + [char]11 + "  3    4" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "8|    |"
$t.Cell(2,2).Range.Text = "14 x 21" + [char]11 + "  2    1" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "4|    |"
$t.Cell(2,3).Range.Text = "71 x 97" + [char]11 + "  9    7" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "1|    |"
$t.Cell(3,1).Range.Text = "97 x 65" + [char]11 + "  6    5" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "7|    |"
$t.Cell(3,2).Range.Text = "67 x 52" + [char]11 + "  5    2" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "7|    |"
$t.Cell(3,3).Range.Text = "78 x 48" + [char]11 + "  4    8" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "8|    |"
$t.Cell(4,1).Range.Text = "23 x 31" + [char]11 + "  3    1" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "3|    |"
$t.Cell(4,2).Range.Text = "22 x 71" + [char]11 + "  7    1" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "2|    |"
$t.Cell(4,3).Range.Text = "47 x 51" + [char]11 + "  5    1" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "7|    |"
$t.Cell(5,1).Range.Text = "22 x 19" + [char]11 + "  1    9" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "2|    |"
$t.Cell(5,2).Range.Text = "86 x 63" + [char]11 + "  6    3" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "6|    |"
$t.Cell(5,3).Range.Text = "60 x 79" + [char]11 + "  7    9" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "0|    |"
